$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.012.66'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '2.741.14'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.26'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.46'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  +5.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.80'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").Value = '3.225.48'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.81'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '63.866.97'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '2.747.24'
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.17'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.81'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.52'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.61'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.521'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -5.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.37'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.43'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").Value = '0.0₃0918'
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.23'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +2.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.35'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +8.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.26'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.93'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.13'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("E35").Value = '  +1.04%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.991'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '349.01'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +5.48%  '
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.61'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -1.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.00'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.14'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0585'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.627'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -1.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.21'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = '2.140.09'
$ws.Range("E51").Value = '  +0.92%  '
